$d = $word.ActiveDocument

$replacement = 'Timeline:^l^l1. May 3, 2024 13:05 - Sajeev Debnath emailed Daniel Parsons asking for the cost of implementing an electrical system in his 2-storey house.^l2. May 3, 2024 13:06 - Daniel Parsons replied, stating their typical charge is around $20k and asked for more details about the rooms and square footage.^l3. May 3, 2024 13:09 - Sajeev Debnath asked for price negotiation due to budget constraints.^l4. May 3, 2024 13:10 - Daniel Parsons asked for confidentiality and inquired about a bathroom on the 2nd floor.^l5. May 3, 2024 13:13 - Sajeev Debnath agreed to keep the price discussions confidential and asked if $14k was possible. Also confirmed the presence of a bathroom on the 2nd floor.^l6. May 3, 2024 13:14 - Daniel Parsons stated that $14k would likely be sufficient, but the final price will depend on labor and materials after a site walkthrough.^l7. May 3, 2024 13:17 - Sajeev Debnath agreed to a site walkthrough on Monday, May 6th at 2pm.'

$d.Content.Find.Execute("This is a test string.", $true, $false, $false, $false, $false,
                         $true, 1, $false, $replacement, 2)
